$wb = $excel.ActiveWorkbook

# ALC  @@ -2273,22 +2273,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1404.65
$ws.Range("I33").Value = 1450.3334
$ws.Range("K33").Value = 1450.3334
$ws.Range("M33").Value = -1221.3334

# ALC  @@ -4079,25 +4079,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4326.3335
$ws.Range("I69").Value = 4490
$ws.Range("J69").Value = 3999
$ws.Range("K69").Value = 13470
$ws.Range("L69").Value = 11997
$ws.Range("M69").Value = -12596
$ws.Range("N69").Value = -13745

# ALC  @@ -4229,25 +4229,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 4326.3335
$ws.Range("I72").Value = 4490
$ws.Range("J72").Value = 3999
$ws.Range("K72").Value = 40410
$ws.Range("L72").Value = 35991
$ws.Range("M72").Value = -36042
$ws.Range("N72").Value = -44727

# ALC  @@ -6457,25 +6457,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2866.7083
$ws.Range("I116").Value = 2780.3333
$ws.Range("J116").Value = 3010.6667
$ws.Range("K116").Value = 2780.3333
$ws.Range("L116").Value = 3010.6667
$ws.Range("M116").Value = 661.6667000000002
$ws.Range("N116").Value = -9894.6667

# ALC  @@ -6702,25 +6702,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2805.5557
$ws.Range("J121").Value = 2881.25
$ws.Range("L121").Value = 8643.75
$ws.Range("N121").Value = -12137.75

# ALC  @@ -7394,22 +7394,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2361.9167
$ws.Range("I135").Value = 2361.9167
$ws.Range("K135").Value = 21257.2503
$ws.Range("M135").Value = -18722.2503

# ALC  @@ -7492,25 +7492,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2988.5925
$ws.Range("I137").Value = 3143.56
$ws.Range("J137").Value = 1051.5
$ws.Range("K137").Value = 9430.68
$ws.Range("L137").Value = 3154.5
$ws.Range("M137").Value = -6880.68
$ws.Range("N137").Value = -8254.5

# ALC  @@ -7544,25 +7544,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 155369.73
$ws.Range("I138").Value = 2398.0557
$ws.Range("J138").Value = 207322.38
$ws.Range("K138").Value = 7194.1671
$ws.Range("L138").Value = 621967.14
$ws.Range("M138").Value = -2054.1671
$ws.Range("N138").Value = -632247.14

# ARM  @@ -10719,25 +10719,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2714.6316
$ws.Range("I61").Value = 1763.8
$ws.Range("J61").Value = 3771.111
$ws.Range("K61").Value = 1763.8
$ws.Range("L61").Value = 3771.111
$ws.Range("M61").Value = -1551.8
$ws.Range("N61").Value = -4195.111

# ARM  @@ -14385,25 +14385,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2714.6316
$ws.Range("I136").Value = 1763.8
$ws.Range("J136").Value = 3771.111
$ws.Range("K136").Value = 5291.4
$ws.Range("L136").Value = 11313.333
$ws.Range("M136").Value = -2741.4
$ws.Range("N136").Value = -16413.333

# CRP  @@ -22722,25 +22722,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 344.27274
$ws.Range("I22").Value = 262.42856
$ws.Range("J22").Value = 487.5
$ws.Range("K22").Value = 262.42856
$ws.Range("L22").Value = 487.5
$ws.Range("M22").Value = 87.57144
$ws.Range("N22").Value = -1187.5

# CRP  @@ -23163,22 +23163,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1856.5424
$ws.Range("I31").Value = 876.8946999999999
$ws.Range("K31").Value = 876.8946999999999
$ws.Range("M31").Value = -581.8946999999999

# CRP  @@ -23310,22 +23310,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1856.5424
$ws.Range("I34").Value = 876.8946999999999
$ws.Range("K34").Value = 876.8946999999999
$ws.Range("M34").Value = -674.8946999999999

# CRP  @@ -28085,25 +28085,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9805876
$ws.Range("I132").Value = 1001.375
$ws.Range("J132").Value = 18521320
$ws.Range("K132").Value = 3004.125
$ws.Range("L132").Value = 55563960
$ws.Range("M132").Value = -474.125
$ws.Range("N132").Value = -55569020

# CRP  @@ -28186,25 +28186,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2112.7058
$ws.Range("I134").Value = 1916.6154
$ws.Range("J134").Value = 2750
$ws.Range("K134").Value = 5749.8462
$ws.Range("L134").Value = 8250
$ws.Range("M134").Value = -3214.8462
$ws.Range("N134").Value = -13320

# CUL  @@ -29556,23 +29556,26 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 4000
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 15000
$ws.Range("N19").Value = -15348

# CUL  @@ -33599,22 +33602,22 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 7194.4443
$ws.Range("I99").Value = 4950
$ws.Range("K99").Value = 14850
$ws.Range("M99").Value = -12604

# CUL  @@ -33942,25 +33945,22 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 8244.362999999999
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 8244.362999999999
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 24733.089
$ws.Range("N106").Value = -26625.089
$ws.Range("M106").ClearContents()

# CUL  @@ -35218,25 +35218,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 952.5
$ws.Range("I131").Value = 250
$ws.Range("J131").Value = 1063.421
$ws.Range("K131").Value = 750
$ws.Range("L131").Value = 3190.263
$ws.Range("M131").Value = 4290
$ws.Range("N131").Value = -13270.263

# GSM  @@ -42233,25 +42233,25 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2059.2812
$ws.Range("I132").Value = 1591.1364
$ws.Range("J132").Value = 3089.2
$ws.Range("K132").Value = 4773.4092
$ws.Range("L132").Value = 9267.599999999999
$ws.Range("M132").Value = -2243.4092
$ws.Range("N132").Value = -14327.6

# LTW  @@ -43062,25 +43062,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 125002330
$ws.Range("I7").Value = 200001380
$ws.Range("J7").Value = 3903.3333
$ws.Range("K7").Value = 200001380
$ws.Range("L7").Value = 3903.3333
$ws.Range("M7").Value = -200001268
$ws.Range("N7").Value = -4127.3333

# LTW  @@ -45711,25 +45711,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3717.84
$ws.Range("I61").Value = 3596.4443
$ws.Range("J61").Value = 4030
$ws.Range("K61").Value = 3596.4443
$ws.Range("L61").Value = 4030
$ws.Range("M61").Value = -3394.4443
$ws.Range("N61").Value = -4434

# LTW  @@ -48265,25 +48265,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3717.84
$ws.Range("I113").Value = 3596.4443
$ws.Range("J113").Value = 4030
$ws.Range("K113").Value = 3596.4443
$ws.Range("L113").Value = 4030
$ws.Range("M113").Value = -1426.4443
$ws.Range("N113").Value = -8370

# LTW  @@ -48890,25 +48890,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 125002330
$ws.Range("I126").Value = 200001380
$ws.Range("J126").Value = 3903.3333
$ws.Range("K126").Value = 600004140
$ws.Range("L126").Value = 11709.9999
$ws.Range("M126").Value = -600001670
$ws.Range("N126").Value = -16649.9999

# LTW  @@ -49377,25 +49377,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 55559120
$ws.Range("J136").Value = 83334184
$ws.Range("L136").Value = 250002552
$ws.Range("N136").Value = -250007652

# WVR  @@ -52821,22 +52821,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 39600
$ws.Range("J64").Value = 39600
$ws.Range("L64").Value = 39600
$ws.Range("N64").Value = -40096

# WVR  @@ -52974,22 +52974,22 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 39600
$ws.Range("J67").Value = 39600
$ws.Range("L67").Value = 39600
$ws.Range("N67").Value = -41316

# WVR  @@ -53605,22 +53605,19 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# WVR  @@ -53752,22 +53749,19 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# WVR  @@ -56120,25 +56114,25 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3403735.5
$ws.Range("I132").Value = 2674.84
$ws.Range("J132").Value = 6946507
$ws.Range("K132").Value = 8024.52
$ws.Range("L132").Value = 20839521
$ws.Range("M132").Value = -5494.52
$ws.Range("N132").Value = -20844581

